$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append a new task row at the end (row 13 before insertion below shifts it to 14):
#    "mejorar formulario de bancos para agregar cuenta bancaria" / "no comenzado"
$ws.Cells.Item(13,1).Value2 = "mejorar formulario de bancos para agregar cuenta bancaria"
$ws.Cells.Item(13,2).Value2 = "no comenzado"

# 2) Insert a new task row at row 11 (pushes existing rows 11-13 down to 12-14):
#    "proveedores pagos, arreglar calculo de retenciones" / "no comenzado"
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11,1).Value2 = "proveedores pagos, arreglar calculo de retenciones"
$ws.Cells.Item(11,2).Value2 = "no comenzado"

# 3) Swap the fill colors used by the "no comenzado" (red) and "terminado" (green)
#    conditional formatting rules, and swap which rule has which priority, while
#    keeping the overall text->color mapping (no comenzado=red, en proceso=yellow,
#    terminado=green) the same - this mirrors reordering the rules in Excel's
#    "Manage Rules" dialog.
$rng = $ws.Range("B1:B1048576")
$fcs = $rng.FormatConditions
$fcTerminado = $fcs.Item(1)
$fcNoComenzado = $fcs.Item(3)

$fcTerminado.Formula1 = '="no comenzado"'
$fcTerminado.Priority = 1

$fcNoComenzado.Formula1 = '="terminado"'
$fcNoComenzado.Priority = 3

# Re-fetch the collection so the color updates below are applied to the
# now-current rule objects and are not discarded.
$fcs2 = $rng.FormatConditions
$fcNowNoComenzado = $fcs2.Item(1)
$fcNowTerminado = $fcs2.Item(3)

$fcNowNoComenzado.Interior.Color = 255     # BGR for RGB FF0000 (red)
$fcNowTerminado.Interior.Color = 5296274   # BGR for RGB 92D050 (green)

# 4) Update the selected cell to reflect where the user ended up (C12)
$ws.Range("C12").Select() | Out-Null
